# paycell_update_history_app_store.xlsx -- "update history new datas found"
# Insert newly-discovered app-version/update rows into the history table and
# correct a couple of mis-recorded dates, matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert rows for newly found versions -------------------------------
# Three rows between the existing "6.7.2" row (46) and "6.5.3" row (47).
$ws.Rows("47:49").Insert()
# Four rows after the existing "6.4.1" row (which is now at 54).
$ws.Rows("55:58").Insert()

# --- 2. Fix existing rows whose recorded date was wrong ---------------------
$ws.Range("A46").Value = "'26-05-2023"
$ws.Range("A53").Value = "'25-01-2023"
$ws.Range("A54").Value = "'22-12-2022"

# --- 3. Fill newly-inserted rows (date, version) ----------------------------
$ws.Range("A47").Value = "'16-05-2023"
$ws.Range("B47").Value = "6.7.1"

$ws.Range("A48").Value = "'3-05-2023"
$ws.Range("B48").Value = "6.7.0"

$ws.Range("A49").Value = "'18-04-2023"
$ws.Range("B49").Value = "6.6.0"

$ws.Range("A55").Value = "'22-11-2022"
$ws.Range("B55").Value = "6.2.0"

$ws.Range("A56").Value = "'23-09-2022"
$ws.Range("B56").Value = "6.0.1"

$ws.Range("A57").Value = "'19-09-2022"
$ws.Range("B57").Value = "6.0.0"

$ws.Range("A58").Value = "'12-08-2022"
$ws.Range("B58").Value = "5.9.0"

# --- 4. Wrap long "Update" text in a handful of rows + widen column C -------
$ws.Columns("C").ColumnWidth = 255.6640625

$ws.Range("C13:C16").WrapText = $true
$ws.Range("C19:C21").WrapText = $true

$ws.Rows(13).RowHeight = 28.8
$ws.Rows(14).RowHeight = 57.6
$ws.Rows(15).RowHeight = 57.6
$ws.Rows(16).RowHeight = 57.6
$ws.Rows(19).RowHeight = 28.8
$ws.Rows(20).RowHeight = 129.6
$ws.Rows(21).RowHeight = 57.6

# --- 5. Restore view state ---------------------------------------------------
$ws.Range("B58").Select()
